$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = 0.5
$ws.Range("K1").Value = 0.6

for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 11).Value = 0.6
}

$ws.Range("K1:K51").Select() | Out-Null
